# Auto-generated edit script: update cryptos list (prices/volumes) per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.833.77"
$ws.Range("E2").Value = "  -1.36%  "

$ws.Range("D3").Value = "3.170.50"
$ws.Range("E3").Value = "  -4.38%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'590.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.14%  "

$ws.Range("D6").Value = "'136.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.76%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "3.168.28"
$ws.Range("E8").Value = "  -4.42%  "

$ws.Range("D9").Value = "'0.512"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.36%  "

$ws.Range("E10").Value = "  -4.82%  "

$ws.Range("D11").Value = "'5.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.80%  "

$ws.Range("D12").Value = "'0.457"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.52%  "

$ws.Range("D13").Value = "'0.0000235"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.02%  "

$ws.Range("D14").Value = "'34.46"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "

$ws.Range("D15").Value = "3.691.92"
$ws.Range("E15").Value = "  -4.37%  "

$ws.Range("E16").Value = "  -2.00%  "

$ws.Range("D17").Value = "3.157.01"
$ws.Range("E17").Value = "  -4.65%  "

$ws.Range("D18").Value = "62.768.84"
$ws.Range("E18").Value = "  -1.58%  "

$ws.Range("D19").Value = "'6.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.87%  "

$ws.Range("D20").Value = "'457.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.58%  "

$ws.Range("D21").Value = "'13.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.49%  "

$ws.Range("D22").Value = "'0.710"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.02%  "

$ws.Range("D23").Value = "'7.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.77%  "

$ws.Range("D24").Value = "'13.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.08%  "

$ws.Range("D25").Value = "'83.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.35%  "

$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "'2.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.15%  "

$ws.Range("B28").Value = "FirstDigitalUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D28").Value = "'0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.22%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.13%  "

$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'6.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.30%  "

$ws.Range("D31").Value = "'2.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.86%  "

$ws.Range("D32").Value = "'27.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.99%  "

$ws.Range("E33").Value = "  -2.53%  "

$ws.Range("D34").Value = "'2.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.99%  "

$ws.Range("E35").Value = "  -6.00%  "

$ws.Range("D36").Value = "'5.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.00%  "

$ws.Range("D37").Value = "'51.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.89%  "

$ws.Range("D38").Value = "0.0₃0713"
$ws.Range("E38").Value = "  -3.74%  "

$ws.Range("D39").Value = "'0.0387"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.28%  "

$ws.Range("D40").Value = "'2.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.94%  "

$ws.Range("D41").Value = "'399.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.72%  "

$ws.Range("D42").Value = "'8.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.90%  "

$ws.Range("D43").Value = "'0.112"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.02%  "

$ws.Range("D44").Value = "2.770.88"
$ws.Range("E44").Value = "  -9.52%  "

$ws.Range("D45").Value = "'0.253"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.00%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.64%  "

$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("D48").Value = "'125.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.20%  "

$ws.Range("D49").Value = "'25.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.07%  "

$ws.Range("D50").Value = "'34.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.58%  "

$ws.Range("E51").Value = "  -2.78%  "
